# Agregue la hoja Noviembre
# Duplicate the "Octubre 2021" sheet (with all its data, styles, column
# widths, conditional formatting and data validations) and place the
# copy right after it, renamed to "Noviembre 2021". The original sheet
# order was: Octubre 2021, Hoja2 -> becomes: Octubre 2021, Noviembre 2021, Hoja2.

$wb = $excel.ActiveWorkbook

$octubre = $wb.Worksheets.Item("Octubre 2021")

# Copy() duplicates the sheet (values, styles, cols, conditional
# formatting, data validations, ...) and inserts the copy right after
# $octubre, making it the new active sheet/tab.
$octubre.Copy($null, $octubre)

$noviembre = $wb.Worksheets.Item(2)
$noviembre.Name = "Noviembre 2021"

# Match the author's selection on the new sheet and keep it the active tab.
$noviembre.Activate()
$noviembre.Range("I18").Select() | Out-Null
